# Auto-generated script applying the Sagittarius_Profits market-data refresh diff.
# For each sheet, update changed cells: set new value, or clear cell if removed in diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 908.125
$ws.Range("J17").Value = 908.125
$ws.Range("L17").Value = 2724.375
$ws.Range("N17").Value = -3060.375
$ws.Range("H80").Value = 2128.1875
$ws.Range("I80").Value = 1285.8572
$ws.Range("K80").Value = 3857.5716
$ws.Range("M80").Value = -2859.5716
$ws.Range("H83").Value = 2128.1875
$ws.Range("I83").Value = 1285.8572
$ws.Range("K83").Value = 11572.7148
$ws.Range("M83").Value = -6580.7148
$ws.Range("H100").Value = 3375
$ws.Range("J100").Value = 3375
$ws.Range("L100").Value = 3375
$ws.Range("N100").Value = -4457
$ws.Range("H106").Value = 1975
$ws.Range("I106").Value = 1975
$ws.Range("K106").Value = 1975
$ws.Range("M106").Value = -1344
$ws.Range("H125").Value = 14906.4
$ws.Range("I125").Value = 1766
$ws.Range("J125").Value = 23666.666
$ws.Range("K125").Value = 15894
$ws.Range("L125").Value = 212999.994
$ws.Range("M125").Value = -13434
$ws.Range("N125").Value = -217919.994
$ws.Range("H132").Value = 2249.75
$ws.Range("I132").Value = 2249.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6749.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4219.25
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 1465.6086
$ws.Range("I137").Value = 1061
$ws.Range("K137").Value = 3183
$ws.Range("M137").Value = -633
$ws.Range("H139").Value = 95000
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1247.6666
$ws.Range("I2").Value = 1702.5
$ws.Range("K2").Value = 1702.5
$ws.Range("M2").Value = -1589.5
$ws.Range("H32").Value = 2486472
$ws.Range("I32").Value = 3503535.2
$ws.Range("J32").Value = 637266.0600000001
$ws.Range("K32").Value = 3503535.2
$ws.Range("L32").Value = 637266.0600000001
$ws.Range("M32").Value = -3503248.2
$ws.Range("N32").Value = -637840.0600000001
$ws.Range("H74").Value = 5999.6
$ws.Range("I74").Value = 5999.6
$ws.Range("K74").Value = 5999.6
$ws.Range("M74").Value = -5125.6
$ws.Range("H77").Value = 5999.6
$ws.Range("I77").Value = 5999.6
$ws.Range("K77").Value = 29998
$ws.Range("M77").Value = -25630
$ws.Range("H116").Value = 1247.6666
$ws.Range("I116").Value = 1702.5
$ws.Range("K116").Value = 1702.5
$ws.Range("M116").Value = 591.5
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1247.6666
$ws.Range("I3").Value = 1702.5
$ws.Range("K3").Value = 1702.5
$ws.Range("M3").Value = -1588.5
$ws.Range("H80").Value = 634.625
$ws.Range("I80").Value = 769.25
$ws.Range("J80").Value = 500
$ws.Range("K80").Value = 769.25
$ws.Range("L80").Value = 500
$ws.Range("M80").Value = 228.75
$ws.Range("N80").Value = -2496
$ws.Range("H83").Value = 634.625
$ws.Range("I83").Value = 769.25
$ws.Range("J83").Value = 500
$ws.Range("K83").Value = 3846.25
$ws.Range("L83").Value = 2500
$ws.Range("M83").Value = 1145.75
$ws.Range("N83").Value = -12484
$ws.Range("H94").Value = 1130.75
$ws.Range("I94").Value = 1071.1428
$ws.Range("K94").Value = 1071.1428
$ws.Range("M94").Value = -620.1428000000001
$ws.Range("H105").Value = 2199.7273
$ws.Range("I105").Value = 2188.5557
$ws.Range("K105").Value = 2188.5557
$ws.Range("M105").Value = -441.5556999999999
$ws.Range("H138").Value = 56663.668
$ws.Range("J138").Value = 56663.668
$ws.Range("L138").Value = 56663.668
$ws.Range("N138").Value = -66943.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2477.1
$ws.Range("I31").Value = 2358.875
$ws.Range("K31").Value = 2358.875
$ws.Range("M31").Value = -2063.875
$ws.Range("H34").Value = 2477.1
$ws.Range("I34").Value = 2358.875
$ws.Range("K34").Value = 2358.875
$ws.Range("M34").Value = -2156.875
$ws.Range("H94").Value = 72998.94
$ws.Range("I94").Value = 125161.89
$ws.Range("J94").Value = 5932.2856
$ws.Range("K94").Value = 125161.89
$ws.Range("L94").Value = 5932.2856
$ws.Range("M94").Value = -124710.89
$ws.Range("N94").Value = -6834.2856
$ws.Range("H105").Value = 2952.2
$ws.Range("I105").Value = 1993
$ws.Range("K105").Value = 1993
$ws.Range("M105").Value = -246

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 237.3077
$ws.Range("J12").Value = 142.16667
$ws.Range("L12").Value = 426.50001
$ws.Range("N12").Value = -772.50001
$ws.Range("H40").Value = 57.714287
$ws.Range("I40").Value = 47.75
$ws.Range("J40").Value = 71
$ws.Range("K40").Value = 191
$ws.Range("L40").Value = 284
$ws.Range("M40").Value = -122
$ws.Range("N40").Value = -422
$ws.Range("H70").Value = 14286.714
$ws.Range("I70").Value = 12
$ws.Range("J70").Value = 16665.834
$ws.Range("K70").Value = 36
$ws.Range("L70").Value = 49997.50199999999
$ws.Range("M70").Value = 279
$ws.Range("N70").Value = -50627.50199999999
$ws.Range("H73").Value = 14286.714
$ws.Range("I73").Value = 12
$ws.Range("J73").Value = 16665.834
$ws.Range("K73").Value = 36
$ws.Range("L73").Value = 49997.50199999999
$ws.Range("M73").Value = 1056
$ws.Range("N73").Value = -52181.50199999999
$ws.Range("H131").Value = 590706.3
$ws.Range("J131").Value = 669331.8
$ws.Range("L131").Value = 2007995.4
$ws.Range("N131").Value = -2018075.4
$ws.Range("H132").Value = 2536.5557
$ws.Range("J132").Value = 2058.3333
$ws.Range("L132").Value = 18524.9997
$ws.Range("N132").Value = -23584.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 466.66666
$ws.Range("I97").Value = 550
$ws.Range("J97").Value = 300
$ws.Range("K97").Value = 550
$ws.Range("L97").Value = 300
$ws.Range("M97").Value = -54
$ws.Range("N97").Value = -1292
$ws.Range("H122").Value = 1296.8889
$ws.Range("I122").Value = 1296.8889
$ws.Range("K122").Value = 3890.6667
$ws.Range("M122").Value = -1440.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2435.1428
$ws.Range("I22").Value = 1709.6666
$ws.Range("J22").Value = 6788
$ws.Range("K22").Value = 1709.6666
$ws.Range("L22").Value = 6788
$ws.Range("M22").Value = -1414.6666
$ws.Range("N22").Value = -7378
$ws.Range("H27").Value = 2435.1428
$ws.Range("I27").Value = 1709.6666
$ws.Range("J27").Value = 6788
$ws.Range("K27").Value = 1709.6666
$ws.Range("L27").Value = 6788
$ws.Range("M27").Value = -1602.6666
$ws.Range("N27").Value = -7002
$ws.Range("H46").Value = 3268.2727
$ws.Range("J46").Value = 3833
$ws.Range("L46").Value = 3833
$ws.Range("N46").Value = -4209
$ws.Range("H55").Value = 1474.762
$ws.Range("J55").Value = 2105
$ws.Range("L55").Value = 2105
$ws.Range("N55").Value = -2451
$ws.Range("H63").Value = 84038
$ws.Range("J63").Value = 83999
$ws.Range("L63").Value = 83999
$ws.Range("N63").Value = -85497
$ws.Range("H66").Value = 84038
$ws.Range("J66").Value = 83999
$ws.Range("L66").Value = 251997
$ws.Range("N66").Value = -259485
$ws.Range("H68").Value = 4000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251
$ws.Range("H71").Value = 4000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256
$ws.Range("H136").Value = 4996.4
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 6491
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 19473
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -24573

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 39749.5
$ws.Range("J15").Value = 39749.5
$ws.Range("L15").Value = 39749.5
$ws.Range("N15").Value = -40325.5
$ws.Range("H41").Value = 19690.572
$ws.Range("I41").Value = 19632.334
$ws.Range("K41").Value = 19632.334
$ws.Range("M41").Value = -19242.334
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H132").Value = 2405.4285
$ws.Range("J132").Value = 1874.6
$ws.Range("L132").Value = 5623.799999999999
$ws.Range("N132").Value = -10683.8
$ws.Range("H136").Value = 2565.1667
$ws.Range("I136").Value = 2660.2856
$ws.Range("K136").Value = 7980.8568
$ws.Range("M136").Value = -5430.8568
